$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 20:16"

# Refresh per-country totals. The underlying feed is sorted by total cases,
# so some rows now correspond to a different country than before (e.g. Israel
# overtook Brasil, Colombia overtook Panama, Lituania overtook Armenia, and
# Costa Rica overtook Kuwait) while keeping the same row position.
$ws.Range("A6").Value = "Estados Unidos"
$ws.Range("B6").Value = 62852
$ws.Range("C6").Value = 7996
$ws.Range("D6").Value = 392
$ws.Range("E6").Value = 61575
$ws.Range("F6").Value = 1382
$ws.Range("G6").Value = 105
$ws.Range("H6").Value = 885
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 37098
$ws.Range("C8").Value = 4107
$ws.Range("D8").Value = 3547
$ws.Range("E8").Value = 33345
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 47
$ws.Range("H8").Value = 206
$ws.Range("A22").Value = "Israel"
$ws.Range("B22").Value = 2369
$ws.Range("C22").Value = 439
$ws.Range("D22").Value = 58
$ws.Range("E22").Value = 2306
$ws.Range("F22").Value = 37
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 5
$ws.Range("A23").Value = "Brasil"
$ws.Range("B23").Value = 2297
$ws.Range("C23").Value = 50
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 2247
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 48
$ws.Range("A51").Value = "Colombia"
$ws.Range("B51").Value = 470
$ws.Range("C51").Value = 92
$ws.Range("D51").Value = 8
$ws.Range("E51").Value = 458
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 4
$ws.Range("A52").Value = "Panama"
$ws.Range("B52").Value = 443
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 434
$ws.Range("F52").Value = 33
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 8
$ws.Range("A53").Value = "Croacia"
$ws.Range("B53").Value = 442
$ws.Range("C53").Value = 60
$ws.Range("D53").Value = 22
$ws.Range("E53").Value = 419
$ws.Range("F53").Value = 6
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 1
$ws.Range("A54").Value = "Egipto"
$ws.Range("B54").Value = 442
$ws.Range("C54").Value = 40
$ws.Range("D54").Value = 93
$ws.Range("E54").Value = 328
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 21
$ws.Range("A55").Value = "Barein"
$ws.Range("B55").Value = 419
$ws.Range("C55").Value = 27
$ws.Range("D55").Value = 177
$ws.Range("E55").Value = 238
$ws.Range("F55").Value = 2
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 4
$ws.Range("A56").Value = "Hong Kong"
$ws.Range("B56").Value = 410
$ws.Range("C56").Value = 23
$ws.Range("D56").Value = 102
$ws.Range("E56").Value = 304
$ws.Range("F56").Value = 4
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 4
$ws.Range("A57").Value = "Mexico"
$ws.Range("B57").Value = 405
$ws.Range("C57").Value = 38
$ws.Range("D57").Value = 4
$ws.Range("E57").Value = 396
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 5
$ws.Range("A58").Value = "Estonia"
$ws.Range("B58").Value = 404
$ws.Range("C58").Value = 35
$ws.Range("D58").Value = 8
$ws.Range("E58").Value = 395
$ws.Range("F58").Value = 5
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1
$ws.Range("A59").Value = "Republica Dominicana"
$ws.Range("B59").Value = 392
$ws.Range("C59").Value = 80
$ws.Range("D59").Value = 3
$ws.Range("E59").Value = 379
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 10
$ws.Range("A60").Value = "Argentina"
$ws.Range("B60").Value = 387
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 52
$ws.Range("E60").Value = 327
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 8
$ws.Range("A61").Value = "Serbia"
$ws.Range("B61").Value = 384
$ws.Range("C61").Value = 81
$ws.Range("D61").Value = 15
$ws.Range("E61").Value = 365
$ws.Range("F61").Value = 21
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 4
$ws.Range("A66").Value = "Lituania"
$ws.Range("B66").Value = 274
$ws.Range("C66").Value = 65
$ws.Range("D66").Value = 1
$ws.Range("E66").Value = 269
$ws.Range("F66").Value = 1
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 4
$ws.Range("A67").Value = "Armenia"
$ws.Range("B67").Value = 265
$ws.Range("C67").Value = 16
$ws.Range("D67").Value = 16
$ws.Range("E67").Value = 249
$ws.Range("F67").Value = 6
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("A75").Value = "Costa Rica"
$ws.Range("B75").Value = 201
$ws.Range("C75").Value = 24
$ws.Range("D75").Value = 2
$ws.Range("E75").Value = 197
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 2
$ws.Range("A76").Value = "Kuwait"
$ws.Range("B76").Value = 195
$ws.Range("C76").Value = 4
$ws.Range("D76").Value = 43
$ws.Range("E76").Value = 152
$ws.Range("F76").Value = 6
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("A77").Value = "Uruguay"
$ws.Range("B77").Value = 189
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 189
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("A78").Value = "Principado de Andorra"
$ws.Range("B78").Value = 188
$ws.Range("C78").Value = 24
$ws.Range("D78").Value = 1
$ws.Range("E78").Value = 186
$ws.Range("F78").Value = 6
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 1
$ws.Range("A79").Value = "San Marino"
$ws.Range("B79").Value = 187
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 4
$ws.Range("E79").Value = 162
$ws.Range("F79").Value = 12
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 21
